$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.014.20'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.305.20'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.76'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.13'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -1.62%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.43'
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.65'
$ws.Range("E12").Value = '  +5.38%  '
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.90'
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.663.92'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.318.33'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.783'
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.895.64'
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0900'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.55'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.87'
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  +1.86%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.02'
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.60'
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.96'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.15'
$ws.Range("E33").Value = '  +6.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -7.56%  '
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0689'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("E41").Value = '  -0.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.998.19'
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.43'
$ws.Range("E43").Value = '  +3.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.11'
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.82'
$ws.Range("E46").Value = '  +3.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.78'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  -2.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.532.13'
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.45'
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.73'
$ws.Range("E51").Value = '  -0.35%  '

Write-Output "Applied crypto price/volume updates"
